$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row of data to the dictionary table (row 52)
$ws.Range("A52").Value = "timing_neut_dx"
$ws.Range("B52").Value = "Timing of convalescent sample relative to SARS-CoV-2 diagnosis in days (NA unless infected and 2-month serum data available)"

# Copy the style (borders/format) from the row above so the new row matches
$ws.Range("A51:B51").Copy()
$ws.Range("A52:B52").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Update the active cell selection as recorded in the saved file
$ws.Range("C3").Select()
